# Synchronized with the demo rep
#
# The template's "Map" block (rows 9-12) is removed from the RVL sheet.
# Deleting these rows shifts everything below them up by four rows
# (old row 13 -> new row 9, ... old row 23 -> new row 19), which is
# exactly what the commit's diff shows: the sheet's used range shrinks
# from A1:H23 to A1:H19 and the trailing "End"/"of Map" block that used
# to live at row 13 becomes the new row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

$ws.Rows("9:12").Delete()
